$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "326.42"
Set-TextValue "E2" "-1.19%"
Set-TextValue "D3" "44.70"
Set-TextValue "E3" "1.92%"
Set-TextValue "D4" "5.548"
Set-TextValue "D5" "0.08081"
Set-TextValue "E5" "-2.74%"
Set-TextValue "D6" "8.705"
Set-TextValue "E6" "-0.85%"
Set-TextValue "D7" "4.334"
Set-TextValue "E7" "-3.65%"
Set-TextValue "D8" "1.901"
Set-TextValue "E8" "-3.04%"
Set-TextValue "E9" "-7.36%"
Set-TextValue "D10" "0.9486"
Set-TextValue "E10" "1.46%"
Set-TextValue "D11" "0.1180"
Set-TextValue "E11" "-5.82%"
Set-TextValue "D12" "0.1894"
Set-TextValue "E12" "-3.07%"
Set-TextValue "E13" "6.95%"
Set-TextValue "D14" "0.04180"
Set-TextValue "E14" "5.13%"
Set-TextValue "D15" "0.1065"
Set-TextValue "E15" "0.17%"
Set-TextValue "D16" "0.001275"
Set-TextValue "E16" "-2.21%"
Set-TextValue "D17" "0.006073"
Set-TextValue "E17" "2.62%"
Set-TextValue "D18" "3.602"
Set-TextValue "E18" "2.35%"
Set-TextValue "D20" "8.357"
Set-TextValue "E20" "-7.58%"
Set-TextValue "E21" "0.08%"
Set-TextValue "E22" "3.54%"
Set-TextValue "D23" "0.04252"
Set-TextValue "E23" "-3.34%"
Set-TextValue "D24" "0.001234"
Set-TextValue "E24" "-1.68%"
Set-TextValue "D25" "0.004603"
Set-TextValue "E25" "4.61%"
Set-TextValue "E26" "3.58%"
Set-TextValue "D27" "0.0003996"
Set-TextValue "E27" "0.10%"
Set-TextValue "D39" "0.02664"
Set-TextValue "E39" "-4.85%"
Set-TextValue "D40" "0.05567"
Set-TextValue "E40" "-0.36%"
Set-TextValue "E41" "24.86%"
Set-TextValue "D42" "0.007691"
Set-TextValue "E42" "-2.74%"
Set-TextValue "E43" "-2.00%"
Set-TextValue "D45" "0.009196"
Set-TextValue "E45" "-11.64%"
Set-TextValue "D46" "0.00007122"
Set-TextValue "E46" "-1.30%"
Set-TextValue "D47" "0.00000000751"
Set-TextValue "E47" "0.13%"
Set-TextValue "D48" "0.003438"
Set-TextValue "E48" "-13.23%"
Set-TextValue "D49" "0.002274"
Set-TextValue "E49" "-0.24%"
Set-TextValue "D50" "0.00002104"
Set-TextValue "E50" "0.13%"
Set-TextValue "D51" "0.0002004"
Set-TextValue "E51" "0.13%"
